# Automated map update (2025-07-25 07:05:42)
# - Row 56 (OT 807150729): Caso corrected, Direccion reformatted, Observaciones updated
# - Row for case 6376 (BOYACA 712, OT 808099366) removed entirely, shifting
#   subsequent rows (old 67-69) up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 56 ---
# "Caso" (A56) must stay a TEXT value ("6494"), not get auto-coerced to a
# number, so force text format, assign, then restore default ("Normal")
# style so we don't leave a stray number-format behind on the cell.
$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = "6494"
$ws.Range("A56").Style = "Normal"

$ws.Range("C56").Value = "SEGUI, JUAN FRANCISCO 4507"
$ws.Range("H56").Value = "Picada"

# --- Remove the BOYACA 712 / case 6376 row (currently row 66) ---
# This shifts rows 67-69 up to 66-68 and shrinks the used range to A1:P68,
# matching the dataset's natural "delete one record" refresh.
$ws.Rows(66).Delete()
